# Weekly data refresh for the "Rabanito" price series.
#
# The underlying source data shifted down by one row starting at row 342
# (a new weekly observation was inserted there), and a second new weekly
# observation was inserted further down (ending up at row 457 after the
# first shift). Everything else keeps its original values, just moved.
#
# We replicate this with two native row-inserts (which push the existing
# rows down and carry formatting along for free), then stamp the new rows'
# cells with the boilerplate from their neighbour (same market / category /
# quality / unit / origin / percentage / classification on every row in
# this block) plus the handful of genuinely new values called out in the
# diff (Fecha, Volumen, Precio minimo/maximo/promedio, Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First new row: inserted before existing row 342 -----------------
$ws.Rows.Item(342).Insert()

# Clone the (now pushed-down) neighbouring row's boilerplate columns,
# then overwrite the columns that actually carry new data for this row.
$ws.Range("A343:R343").Copy()
$ws.Range("A342").PasteSpecial()

$ws.Range("D342").Value = 45120
$ws.Range("J342").Value = 7000
$ws.Range("K342").Value = 3000
$ws.Range("L342").Value = 4000
$ws.Range("M342").Value = 3500
$ws.Range("P342").Value = 35

# --- Second new row: inserted before (post-shift) row 457 ------------
$ws.Rows.Item(457).Insert()

$ws.Range("A458:R458").Copy()
$ws.Range("A457").PasteSpecial()

$ws.Range("D457").Value = 45121
$ws.Range("J457").Value = 7000
$ws.Range("K457").Value = 3000
$ws.Range("L457").Value = 4000
$ws.Range("M457").Value = 3500
$ws.Range("P457").Value = 35
